$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking price values to remain text (matching source formatting)
$textCells = $excel.Union($ws.Range("D5"), $ws.Range("D6"), $ws.Range("D9"), $ws.Range("D10"), $ws.Range("D11"), $ws.Range("D13"), $ws.Range("D14"), $ws.Range("D17"), $ws.Range("D19"), $ws.Range("D23"), $ws.Range("D24"), $ws.Range("D27"), $ws.Range("D30"), $ws.Range("D31"), $ws.Range("D32"), $ws.Range("D34"), $ws.Range("D35"), $ws.Range("D36"), $ws.Range("D37"), $ws.Range("D39"), $ws.Range("D40"), $ws.Range("D41"), $ws.Range("D42"), $ws.Range("D43"), $ws.Range("D47"), $ws.Range("D48"), $ws.Range("D51"))
foreach ($area in $textCells.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "49.551.34"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "2.641.28"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "112.74"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").Value = "326.34"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.549"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").Value = "39.67"
$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").Value = "19.99"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("D13").Value = "0.128"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").Value = "7.57"
$ws.Range("E14").Value = "  +2.79%  "

$ws.Range("D15").Value = "3.055.27"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "2.640.31"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").Value = "0.858"
$ws.Range("E17").Value = "  -1.91%  "

$ws.Range("D18").Value = "49.495.01"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "13.35"
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").Value = "0.0₃0949"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").Value = "268.10"
$ws.Range("E23").Value = "  -3.33%  "

$ws.Range("D24").Value = "69.17"
$ws.Range("E24").Value = "  -4.10%  "

$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "26.06"

$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").Value = "0.137"
$ws.Range("E30").Value = "  -2.83%  "

$ws.Range("D31").Value = "34.70"
$ws.Range("E31").Value = "  -3.87%  "

$ws.Range("D32").Value = "49.64"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").Value = "0.0820"
$ws.Range("E34").Value = "  +0.94%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "19.11"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("D37").Value = "4.93"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "23.88"
$ws.Range("E40").Value = "  +8.11%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "129.47"
$ws.Range("E41").Value = "  +4.41%  "

$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  +3.36%  "

$ws.Range("D43").Value = "0.0341"
$ws.Range("E43").Value = "  +8.36%  "

$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").Value = "2.061.40"
$ws.Range("E45").Value = "  -1.25%  "

$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +6.31%  "

$ws.Range("D48").Value = "2.21"
$ws.Range("E48").Value = "  -6.11%  "

$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("E50").Value = "  -3.02%  "

$ws.Range("D51").Value = "58.87"
$ws.Range("E51").Value = "  -1.10%  "
